# Auto-generated edit script applying numeric corrections to the Leve profit sheets
# per the commit diff (scheduled runner update of cached market-board prices).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 333
$ws.Range("I33").Value = 347.25
$ws.Range("K33").Value = 347.25
$ws.Range("M33").Value = -118.25

$ws.Range("H103").Value = 528.1667
$ws.Range("I103").Value = 481.8
$ws.Range("K103").Value = 1445.4
$ws.Range("M103").Value = -859.4000000000001

$ws.Range("H113").Value = 5271.1816
$ws.Range("I113").Value = 3248
$ws.Range("K113").Value = 3248
$ws.Range("M113").Value = 6

$ws.Range("H123").Value = 45499
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 45499
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 45499
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -55299

$ws.Range("H138").Value = 1253.5385
$ws.Range("I138").Value = 941.3333
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 2823.9999
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 2316.0001
$ws.Range("N138").Value = -25280


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 12689.333
$ws.Range("I37").Value = 12689.333
$ws.Range("K37").Value = 12689.333
$ws.Range("M37").Value = -12416.333

$ws.Range("H43").Value = 29962.334
$ws.Range("J43").Value = 29962.334
$ws.Range("L43").Value = 29962.334
$ws.Range("N43").Value = -30588.334

$ws.Range("H63").Value = 11262.235
$ws.Range("I63").Value = 11243
$ws.Range("J63").Value = 11324.75
$ws.Range("K63").Value = 11243
$ws.Range("L63").Value = 11324.75
$ws.Range("M63").Value = -10557
$ws.Range("N63").Value = -12696.75

$ws.Range("H66").Value = 11262.235
$ws.Range("I66").Value = 11243
$ws.Range("J66").Value = 11324.75
$ws.Range("K66").Value = 56215
$ws.Range("L66").Value = 56623.75
$ws.Range("M66").Value = -52783
$ws.Range("N66").Value = -63487.75

$ws.Range("H74").Value = 1741.3334
$ws.Range("I74").Value = 1680.6
$ws.Range("J74").Value = 1893.1666
$ws.Range("K74").Value = 1680.6
$ws.Range("L74").Value = 1893.1666
$ws.Range("M74").Value = -806.5999999999999
$ws.Range("N74").Value = -3641.1666

$ws.Range("H77").Value = 1741.3334
$ws.Range("I77").Value = 1680.6
$ws.Range("J77").Value = 1893.1666
$ws.Range("K77").Value = 8403
$ws.Range("L77").Value = 9465.833000000001
$ws.Range("M77").Value = -4035
$ws.Range("N77").Value = -18201.833

$ws.Range("H80").Value = 85209.53
$ws.Range("J80").Value = 85209.53
$ws.Range("L80").Value = 85209.53
$ws.Range("N80").Value = -87205.53

$ws.Range("H83").Value = 85209.53
$ws.Range("J83").Value = 85209.53
$ws.Range("L83").Value = 255628.59
$ws.Range("N83").Value = -265612.59


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1927.2222
$ws.Range("I20").Value = 2034.5714
$ws.Range("K20").Value = 2034.5714
$ws.Range("M20").Value = -1787.5714

$ws.Range("H22").Value = 359.5
$ws.Range("I22").Value = 456
$ws.Range("K22").Value = 456
$ws.Range("M22").Value = -283

$ws.Range("H35").Value = 48000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws.Range("H134").Value = 4789.1562
$ws.Range("I134").Value = 5043.375
$ws.Range("J134").Value = 4026.5
$ws.Range("K134").Value = 15130.125
$ws.Range("L134").Value = 12079.5
$ws.Range("M134").Value = -12595.125
$ws.Range("N134").Value = -17149.5


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H19").Value = 12924.625
$ws.Range("I19").Value = 16966.166
$ws.Range("J19").Value = 800
$ws.Range("K19").Value = 16966.166
$ws.Range("L19").Value = 800
$ws.Range("M19").Value = -16796.166
$ws.Range("N19").Value = -1140

$ws.Range("H24").Value = 12924.625
$ws.Range("I24").Value = 16966.166
$ws.Range("J24").Value = 800
$ws.Range("K24").Value = 16966.166
$ws.Range("L24").Value = 800
$ws.Range("M24").Value = -16796.166
$ws.Range("N24").Value = -1140

$ws.Range("H94").Value = 1316.3334
$ws.Range("I94").Value = 1316.3334
$ws.Range("K94").Value = 1316.3334
$ws.Range("M94").Value = -865.3334

$ws.Range("H107").Value = 323.66666
$ws.Range("I107").Value = 262.18182
$ws.Range("K107").Value = 262.18182
$ws.Range("M107").Value = 1657.81818


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2804223.2
$ws.Range("I4").Value = 93661.664
$ws.Range("K4").Value = 280984.992
$ws.Range("M4").Value = -280872.992

$ws.Range("H92").Value = 307.25
$ws.Range("I92").Value = 307.25
$ws.Range("K92").Value = 921.75
$ws.Range("M92").Value = 326.25

$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

$ws.Range("H131").Value = 1357

$ws.Range("H137").Value = 6139.6665
$ws.Range("I137").Value = 3970
$ws.Range("J137").Value = 7224.5
$ws.Range("K137").Value = 11910
$ws.Range("L137").Value = 21673.5
$ws.Range("M137").Value = -6810
$ws.Range("N137").Value = -31873.5

$ws.Range("H140").Value = 627050
$ws.Range("I140").Value = 627050
$ws.Range("K140").Value = 1881150
$ws.Range("M140").Value = -1875970


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 14645
$ws.Range("I57").Value = 15000
$ws.Range("J57").Value = 14290
$ws.Range("K57").Value = 15000
$ws.Range("L57").Value = 14290
$ws.Range("M57").Value = -14180
$ws.Range("N57").Value = -15930

$ws.Range("H132").Value = 4877.154
$ws.Range("I132").Value = 4521.2
$ws.Range("K132").Value = 13563.6
$ws.Range("M132").Value = -11033.6


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 955.7
$ws.Range("I22").Value = 755.4
$ws.Range("J22").Value = 1156
$ws.Range("K22").Value = 755.4
$ws.Range("L22").Value = 1156
$ws.Range("M22").Value = -460.4
$ws.Range("N22").Value = -1746

$ws.Range("H27").Value = 955.7
$ws.Range("I27").Value = 755.4
$ws.Range("J27").Value = 1156
$ws.Range("K27").Value = 755.4
$ws.Range("L27").Value = 1156
$ws.Range("M27").Value = -648.4
$ws.Range("N27").Value = -1370

$ws.Range("H33").Value = 10526
$ws.Range("J33").Value = 10526
$ws.Range("L33").Value = 10526
$ws.Range("N33").Value = -11106

$ws.Range("H94").Value = 42332.668
$ws.Range("J94").Value = 42332.668
$ws.Range("L94").Value = 42332.668
$ws.Range("N94").Value = -43684.668

$ws.Range("H98").Value = 58999.5
$ws.Range("J98").Value = 58999.5
$ws.Range("L98").Value = 58999.5
$ws.Range("N98").Value = -64989.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H62").Value = 8334
$ws.Range("I62").Value = 6001
$ws.Range("J62").Value = 13000
$ws.Range("K62").Value = 6001
$ws.Range("L62").Value = 13000
$ws.Range("M62").Value = -5377
$ws.Range("N62").Value = -14248

$ws.Range("H65").Value = 8334
$ws.Range("I65").Value = 6001
$ws.Range("J65").Value = 13000
$ws.Range("K65").Value = 30005
$ws.Range("L65").Value = 65000
$ws.Range("M65").Value = -26885
$ws.Range("N65").Value = -71240

$ws.Range("H132").Value = 3967.516
$ws.Range("I132").Value = 3606.75
$ws.Range("J132").Value = 7334.6665
$ws.Range("K132").Value = 10820.25
$ws.Range("L132").Value = 22003.9995
$ws.Range("M132").Value = -8290.25
$ws.Range("N132").Value = -27063.9995

